# Apply the change: add a date value to cell C1 on the "About" sheet,
# formatted as a date (m/d/yyyy) -> serial 44307 (2021-04-21).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")
$ws.Range("C1").Value = 44307
$ws.Range("C1").NumberFormat = "m/d/yyyy"
